{"js": "// In-Class Assignment Day 4 - correction:\n//   The \"Book\" relation definition carried \"author\" and \"genre\" as its\n//   own attributes, but the class had already normalized those out into\n//   separate Author / Genre / Book_Genre relations below it. Drop\n//   \", author, genre\" so \"Book\" only keeps \", title\".\n//   Word's hidden \"_GoBack\" bookmark tracks the most recent edit point;\n//   it was sitting at the end of the old \"Book_Genre([ISBN, genre])\"\n//   line and needs to move to right after \"title\" in the \"Book\" line.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the specific \"Book ( [ISBN], title, author, genre) \" paragraph by\n// its exact text, so the many other \"title\"/\"author\"/\"genre\" example\n// mentions earlier in the document are left untouched.\nlet bookParagraph = null;\nfor (const p of paragraphs.items) {\n  if (p.text === \"Book ([ISBN], title, author, genre) \") {\n    bookParagraph = p;\n    break;\n  }\n}\n\nif (bookParagraph) {\n  // Remove \", author, genre\" so only \", title\" remains before the \")\".\n  const removeResults = bookParagraph.getRange().search(\", author, genre\", { matchCase: true });\n  removeResults.load(\"items\");\n  await context.sync();\n\n  if (removeResults.items.length > 0) {\n    removeResults.items[0].insertText(\"\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n\n  // Drop the old \"_GoBack\" bookmark (Word keeps only one instance of it\n  // document-wide) from wherever it used to sit.\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n\n  // Re-locate \"title\" inside the (now shorter) \"Book\" paragraph and put\n  // a collapsed \"_GoBack\" bookmark right after it, before the \")\".\n  const titleResults = bookParagraph.getRange().search(\"title\", { matchCase: true });\n  titleResults.load(\"items\");\n  await context.sync();\n\n  if (titleResults.items.length > 0) {\n    const afterTitle = titleResults.items[0].getRange(Word.RangeLocation.end);\n    afterTitle.insertBookmark(\"_GoBack\");\n    await context.sync();\n  }\n}\n", "ps1": "# In-Class Assignment Day 4 - correction:\n#   \"Book\" relation definition loses the \"author\" and \"genre\" attributes\n#   (they don't belong there once Author/Genre/Book_Genre were normalized\n#   out), leaving just \"title\". The \"_GoBack\" bookmark - which Word had\n#   left at the end of the old \"Book_Genre([ISBN, genre])\" line from the\n#   last edit - ends up at the new last-edited spot, right after \"title\"\n#   in the \"Book\" line.\n\n$wdCollapseEnd = 0\n\n$d = $word.ActiveDocument\n\n# Locate the \"Book ( [ISBN], title, author, genre) \" definition paragraph\n# (scan paragraphs rather than a global Find so we don't touch the many\n# other \"title\"/\"author\"/\"genre\" occurrences used as examples earlier in\n# the document).\n$bookPara = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -match \"^Book \\(\\[ISBN\\], title, author, genre\\)\") {\n        $bookPara = $p\n        break\n    }\n}\n\nif ($bookPara -ne $null) {\n    # Remove \", author, genre\" so only \", title\" remains before the \")\".\n    $editRange = $bookPara.Range\n    $editFind = $editRange.Find\n    $editFind.Text = \", author, genre\"\n    if ($editFind.Execute()) {\n        $editRange.Text = \"\"\n    }\n\n    # Re-use the (now shorter) paragraph range and find \"title\" again so\n    # we can drop a collapsed bookmark right after it - this is where\n    # \"_GoBack\" (Word only ever keeps a single instance of it) now belongs.\n    $bookRange2 = $bookPara.Range\n    $titleFind = $bookRange2.Find\n    $titleFind.Text = \"title\"\n    if ($titleFind.Execute()) {\n        $bookRange2.Collapse($wdCollapseEnd)\n        $d.Bookmarks.Add(\"_GoBack\", $bookRange2)\n    }\n}\n"}
